$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Sheet1"

# Update the allocation table:
#  - drop the "Balite" row entirely
#  - rename "Buguion" -> "Bulusan"
#  - reassign shelters for the remaining communities
#  - drop the trailing "Poblacion" / "San Miguel Meysulao High School" row's old shelter,
#    folding Poblacion in with the rest using "Gatbuca Basketball Court"
#  - rename "San Miguel Meysulao High School" -> "Doña Damiana Elem School" (now assigned to Bulusan)

$ws.Range("A2").Value = "Balungao"
$ws.Range("B2").Value = "Gatbuca Basketball Court"
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = "Bulusan"
$ws.Range("B3").Value = "Doña Damiana Elem School"
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = "Calizon"
$ws.Range("B4").Value = "Gatbuca Basketball Court"
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "Poblacion"
$ws.Range("B5").Value = "Gatbuca Basketball Court"
$ws.Range("C5").Value = 1

# Remove the now-unused 6th row
$ws.Range("A6:C6").Delete()
